# Swap the deck's colour theme from the custom "Integral" / "Red Violet"
# scheme over to the built-in Office Theme ("Office") colour values.
#
# The underlying OOXML change is a full swap of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml (the slide-master theme and the notes-master theme
# trade places). The PowerPoint object model only exposes the slide
# master's theme colours for editing (Slide.ThemeColorScheme /
# Master.ColorScheme), so we drive the 12 theme colour slots to the
# "Office Theme" values through that supported surface.

$p = $ppt.ActivePresentation

# Theme colour slots, in Office "ThemeColorScheme" order:
#  1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#  8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
$officeThemeRGB = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

function ToComRGB([int]$hexRGB) {
    # PowerPoint's RGB longs are stored BGR (0x00BBGGRR), so flip the
    # standard 0xRRGGBB hex value before assigning it.
    $r = ($hexRGB -shr 16) -band 0xFF
    $g = ($hexRGB -shr 8) -band 0xFF
    $b = $hexRGB -band 0xFF
    return ($b -shl 16) -bor ($g -shl 8) -bor $r
}

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeThemeRGB.Length; $i++) {
    $tcs.Item($i).RGB = ToComRGB($officeThemeRGB[$i - 1])
}
